# Apply the "Add files via upload" edit:
#   1. Remove the last two slides ("Pirma programa" / id 269 and
#      "Uzduotis" / id 271), which also drops their notes pages.
#   2. Fix the course repo URL on the "Mokymosi medziaga" slide from
#      .../ca-java-lessons to .../code-academy-java-kursai.

$p = $ppt.ActivePresentation

# --- 1. Delete the two trailing slides -------------------------------
# Match by SlideID (stable even while we delete/reindex) rather than a
# fixed positional index.
$idsToDelete = @(269, 271)
foreach ($targetId in $idsToDelete) {
    for ($i = $p.Slides.Count; $i -ge 1; $i--) {
        $slide = $p.Slides.Item($i)
        if ($slide.SlideID -eq $targetId) {
            $slide.Delete()
            break
        }
    }
}

# --- 2. Update the GitHub URL text -----------------------------------
$oldUrl = "https://github.com/tadus21/ca-java-lessons"
$newUrl = "https://github.com/tadus21/code-academy-java-kursai"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            $fullText = $tr.Text
            $idx = $fullText.IndexOf($oldUrl)
            if ($idx -ge 0) {
                # Replace the whole run's text in one shot (rather than
                # just the differing suffix) so the run stays a single
                # <a:r> element with its original formatting intact.
                $sub = $tr.Characters($idx + 1, $oldUrl.Length)
                $sub.Text = $newUrl
            }
        }
    }
}
